$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06447966666666667
$ws.Range("H2").Value = 0.193439
$ws.Range("I2").Value = 0.001101138907643723
$ws.Range("J2").Value = 0.001101138907643722
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.083188
$ws.Range("N2").Value = 6.249564
$ws.Range("O2").Value = 0.01853451022101116
$ws.Range("P2").Value = 0.01853451022101116
$ws.Range("Q2").Value = 0.134323267844
$ws.Range("R2").Value = 1.208909410596
$ws.Range("S2").Value = 0.00002040907033847564
$ws.Range("T2").Value = 0.00002040907033847564
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06447966666666667
$ws.Range("H3").Value = 0.193439
$ws.Range("I3").Value = 0.001101138907643723
$ws.Range("J3").Value = 0.001101138907643722
$ws.Range("O3").Value = 0.7177032719746937
$ws.Range("P3").Value = 0.717703271974694
$ws.Range("Q3").Value = 5.201337811704667
$ws.Range("R3").Value = 46.812040305342
$ws.Range("S3").Value = 0.0007902909969145398
$ws.Range("T3").Value = 0.0007902909969145399
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06447966666666667
$ws.Range("H4").Value = 0.193439
$ws.Range("I4").Value = 0.001101138907643723
$ws.Range("J4").Value = 0.001101138907643722
$ws.Range("M4").Value = 29.09185666666666
$ws.Range("N4").Value = 87.27556999999999
$ws.Range("O4").Value = 0.258835647448298
$ws.Range("P4").Value = 0.258835647448298
$ws.Range("Q4").Value = 1.875833220581111
$ws.Range("R4").Value = 16.88249898523
$ws.Range("S4").Value = 0.0002850140020904745
$ws.Range("T4").Value = 0.0002850140020904745
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06447966666666667
$ws.Range("H5").Value = 0.193439
$ws.Range("I5").Value = 0.001101138907643723
$ws.Range("J5").Value = 0.001101138907643722
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5537223333333333
$ws.Range("N5").Value = 1.661167
$ws.Range("O5").Value = 0.004926570355997066
$ws.Range("P5").Value = 0.004926570355997067
$ws.Range("Q5").Value = 0.03570383147922222
$ws.Range("R5").Value = 0.321334483313
$ws.Range("S5").Value = 0.000005424838300232554
$ws.Range("T5").Value = 0.000005424838300232554
$ws.Range("I6").Value = 0.00657695954769643
$ws.Range("J6").Value = 0.006576959547696431
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.083188
$ws.Range("N6").Value = 6.249564
$ws.Range("O6").Value = 0.01853451022101116
$ws.Range("P6").Value = 0.01853451022101116
$ws.Range("Q6").Value = 0.8022954168560001
$ws.Range("R6").Value = 7.220658751704001
$ws.Range("S6").Value = 0.0001219007239599564
$ws.Range("T6").Value = 0.0001219007239599565
$ws.Range("I7").Value = 0.00657695954769643
$ws.Range("J7").Value = 0.006576959547696431
$ws.Range("O7").Value = 0.7177032719746937
$ws.Range("P7").Value = 0.717703271974694
$ws.Range("S7").Value = 0.00472030538702693
$ws.Range("T7").Value = 0.004720305387026931
$ws.Range("I8").Value = 0.00657695954769643
$ws.Range("J8").Value = 0.006576959547696431
$ws.Range("M8").Value = 29.09185666666666
$ws.Range("N8").Value = 87.27556999999999
$ws.Range("O8").Value = 0.258835647448298
$ws.Range("P8").Value = 0.258835647448298
$ws.Range("Q8").Value = 11.20410796889111
$ws.Range("R8").Value = 100.83697172002
$ws.Range("S8").Value = 0.00170235158276927
$ws.Range("T8").Value = 0.001702351582769271
$ws.Range("I9").Value = 0.00657695954769643
$ws.Range("J9").Value = 0.006576959547696431
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5537223333333333
$ws.Range("N9").Value = 1.661167
$ws.Range("O9").Value = 0.004926570355997066
$ws.Range("P9").Value = 0.004926570355997067
$ws.Range("Q9").Value = 0.2132543439402222
$ws.Range("R9").Value = 1.919289095462
$ws.Range("S9").Value = 0.0000324018539402731
$ws.Range("T9").Value = 0.00003240185394027311
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008175
$ws.Range("H10").Value = 0.024525
$ws.Range("I10").Value = 0.0001396069650378791
$ws.Range("J10").Value = 0.0001396069650378791
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.083188
$ws.Range("N10").Value = 6.249564
$ws.Range("O10").Value = 0.01853451022101116
$ws.Range("P10").Value = 0.01853451022101116
$ws.Range("Q10").Value = 0.0170300619
$ws.Range("R10").Value = 0.1532705571
$ws.Range("S10").Value = 0.000002587546720418918
$ws.Range("T10").Value = 0.000002587546720418918
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.008175
$ws.Range("H11").Value = 0.024525
$ws.Range("I11").Value = 0.0001396069650378791
$ws.Range("J11").Value = 0.0001396069650378791
$ws.Range("O11").Value = 0.7177032719746937
$ws.Range("P11").Value = 0.717703271974694
$ws.Range("Q11").Value = 0.65944721505
$ws.Range("R11").Value = 5.93502493545
$ws.Range("S11").Value = 0.0001001963755981425
$ws.Range("T11").Value = 0.0001001963755981425
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.008175
$ws.Range("H12").Value = 0.024525
$ws.Range("I12").Value = 0.0001396069650378791
$ws.Range("J12").Value = 0.0001396069650378791
$ws.Range("M12").Value = 29.09185666666666
$ws.Range("N12").Value = 87.27556999999999
$ws.Range("O12").Value = 0.258835647448298
$ws.Range("P12").Value = 0.258835647448298
$ws.Range("Q12").Value = 0.23782592825
$ws.Range("R12").Value = 2.14043335425
$ws.Range("S12").Value = 0.00003613525918387133
$ws.Range("T12").Value = 0.00003613525918387134
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.008175
$ws.Range("H13").Value = 0.024525
$ws.Range("I13").Value = 0.0001396069650378791
$ws.Range("J13").Value = 0.0001396069650378791
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.5537223333333333
$ws.Range("N13").Value = 1.661167
$ws.Range("O13").Value = 0.004926570355997066
$ws.Range("P13").Value = 0.004926570355997067
$ws.Range("Q13").Value = 0.004526680075
$ws.Range("R13").Value = 0.040740120675
$ws.Range("S13").Value = 0.0000006877835354463338
$ws.Range("T13").Value = 0.0000006877835354463341
$ws.Range("G14").Value = 58.099467
$ws.Range("H14").Value = 174.298401
$ws.Range("I14").Value = 0.992182294579622
$ws.Range("J14").Value = 0.992182294579622
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.083188
$ws.Range("N14").Value = 6.249564
$ws.Range("O14").Value = 0.01853451022101116
$ws.Range("P14").Value = 0.01853451022101116
$ws.Range("Q14").Value = 121.032112460796
$ws.Range("R14").Value = 1089.289012147164
$ws.Range("S14").Value = 0.01838961287999231
$ws.Range("T14").Value = 0.01838961287999231
$ws.Range("G15").Value = 58.099467
$ws.Range("H15").Value = 174.298401
$ws.Range("I15").Value = 0.992182294579622
$ws.Range("J15").Value = 0.992182294579622
$ws.Range("O15").Value = 0.7177032719746937
$ws.Range("P15").Value = 0.717703271974694
$ws.Range("Q15").Value = 4686.670545448243
$ws.Range("R15").Value = 42180.03490903418
$ws.Range("S15").Value = 0.7120924792151542
$ws.Range("T15").Value = 0.7120924792151544
$ws.Range("G16").Value = 58.099467
$ws.Range("H16").Value = 174.298401
$ws.Range("I16").Value = 0.992182294579622
$ws.Range("J16").Value = 0.992182294579622
$ws.Range("M16").Value = 29.09185666666666
$ws.Range("N16").Value = 87.27556999999999
$ws.Range("O16").Value = 0.258835647448298
$ws.Range("P16").Value = 0.258835647448298
$ws.Range("Q16").Value = 1690.22136637373
$ws.Range("R16").Value = 15211.99229736357
$ws.Range("S16").Value = 0.2568121466042543
$ws.Range("T16").Value = 0.2568121466042544
$ws.Range("G17").Value = 58.099467
$ws.Range("H17").Value = 174.298401
$ws.Range("I17").Value = 0.992182294579622
$ws.Range("J17").Value = 0.992182294579622
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5537223333333333
$ws.Range("N17").Value = 1.661167
$ws.Range("O17").Value = 0.004926570355997066
$ws.Range("P17").Value = 0.004926570355997067
$ws.Range("Q17").Value = 32.170972432663
$ws.Range("R17").Value = 289.538751893967
$ws.Range("S17").Value = 0.004888055880221114
$ws.Range("T17").Value = 0.004888055880221115
